function Set-CellText($sheet, $addr, $text) {
    $cell = $sheet.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-CellText $ws "D2" "59.736.41"
Set-CellText $ws "E2" "  -1.40%  "
Set-CellText $ws "D3" "2.369.70"
Set-CellText $ws "E3" "  -1.11%  "
Set-CellText $ws "E4" "  +0.37%  "
Set-CellText $ws "D5" "558.79"
Set-CellText $ws "E5" "  -1.84%  "
Set-CellText $ws "D6" "137.50"
Set-CellText $ws "E6" "  -1.56%  "
Set-CellText $ws "E7" "  -0.16%  "
Set-CellText $ws "E8" "  +1.05%  "
Set-CellText $ws "D9" "2.365.16"
Set-CellText $ws "E9" "  -0.43%  "
Set-CellText $ws "E10" "  -2.07%  "
Set-CellText $ws "E11" "  -0.84%  "
Set-CellText $ws "E12" "  +0.80%  "
Set-CellText $ws "E13" "  +0.69%  "
Set-CellText $ws "D14" "25.59"
Set-CellText $ws "E14" "  -1.14%  "
Set-CellText $ws "D15" "2.799.16"
Set-CellText $ws "E15" "  -1.02%  "
Set-CellText $ws "D16" "0.0000164"
Set-CellText $ws "E16" "  -2.86%  "
Set-CellText $ws "D17" "59.670.82"
Set-CellText $ws "E17" "  -1.61%  "
Set-CellText $ws "D18" "2.364.99"
Set-CellText $ws "E18" "  -0.07%  "
Set-CellText $ws "D19" "8.03"
Set-CellText $ws "E19" "  +13.11%  "
Set-CellText $ws "D20" "10.50"
Set-CellText $ws "E20" "  -0.07%  "
Set-CellText $ws "D21" "320.59"
Set-CellText $ws "E21" "  +0.05%  "
Set-CellText $ws "E22" "  +1.09%  "
Set-CellText $ws "D23" "6.05"
Set-CellText $ws "E23" "  -0.77%  "
Set-CellText $ws "E24" "  -0.09%  "
Set-CellText $ws "E25" "  -3.12%  "
Set-CellText $ws "D26" "64.03"
Set-CellText $ws "E26" "  -0.64%  "
Set-CellText $ws "D27" "556.92"
Set-CellText $ws "E27" "  -3.15%  "
Set-CellText $ws "D28" "8.12"
Set-CellText $ws "E28" "  -6.45%  "
Set-CellText $ws "D29" "2.483.51"
Set-CellText $ws "E29" "  -1.24%  "
Set-CellText $ws "D30" "0.0₃0918"
Set-CellText $ws "E30" "  +1.52%  "
Set-CellText $ws "D31" "7.99"
Set-CellText $ws "E31" "  +2.55%  "
Set-CellText $ws "E32" "  -2.55%  "
Set-CellText $ws "E33" "  -2.71%  "
Set-CellText $ws "E34" "  -0.97%  "
Set-CellText $ws "E35" "  -0.49%  "
Set-CellText $ws "E36" "  +3.77%  "
Set-CellText $ws "D37" "152.85"
Set-CellText $ws "E37" "  +3.66%  "
Set-CellText $ws "E38" "  +0.70%  "
Set-CellText $ws "E39" "  -1.46%  "
Set-CellText $ws "D40" "18.10"
Set-CellText $ws "E40" "  +0.29%  "
Set-CellText $ws "D41" "4.97"
Set-CellText $ws "E41" "  -1.69%  "
Set-CellText $ws "E42" "  -0.11%  "
Set-CellText $ws "D43" "41.49"
Set-CellText $ws "E43" "  +0.15%  "
Set-CellText $ws "E44" "  -1.19%  "
Set-CellText $ws "D45" "2.42"
Set-CellText $ws "E45" "  +3.60%  "
Set-CellText $ws "E46" "  +4.19%  "
Set-CellText $ws "D47" "139.00"
Set-CellText $ws "E47" "  -0.80%  "
Set-CellText $ws "E48" "  +0.72%  "
Set-CellText $ws "D49" "0.584"
Set-CellText $ws "E49" "  +0.33%  "
Set-CellText $ws "D50" "0.0498"
Set-CellText $ws "E50" "  -0.70%  "
Set-CellText $ws "D51" "18.97"
Set-CellText $ws "E51" "  -1.54%  "
